# Update the three test e-mail addresses on the "AccountCreationData" sheet
# (newtest4/5/6@gmail.com -> aewtest1/2/3@gmail.com) and move the cursor
# selection to E17, matching the author's latest edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

$ws.Activate()

$ws.Range("A2").Value = "aewtest1@gmail.com"
$ws.Range("A3").Value = "aewtest2@gmail.com"
$ws.Range("A4").Value = "aewtest3@gmail.com"

$ws.Range("E17").Select() | Out-Null
